$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the missing X6 / Y6 values on the existing last row
$ws.Range("X6").Value = -0.35999999999999943
$ws.Range("Y6").Value = "Down"

# Carry the existing date / percentage cell formatting down to row 7
# (reuses the workbook's existing cell styles instead of creating new ones)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("S6:T6").Copy()
$ws.Range("S7:T7").PasteSpecial(-4122)

# Append new row 7 with the latest scan data
$ws.Range("A7").Value = 42648.885277777779
$ws.Range("B7").Value = -8
$ws.Range("C7").Value = "Sell"
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 21108
$ws.Range("F7").Value = 1173
$ws.Range("G7").Value = 53
$ws.Range("H7").Value = 43
$ws.Range("I7").Value = 63
$ws.Range("J7").Value = 36
$ws.Range("K7").Value = 29385
$ws.Range("L7").Value = 190
$ws.Range("M7").Value = 155
$ws.Range("N7").Value = 37
$ws.Range("O7").Value = 21
$ws.Range("P7").Value = "Named"
$ws.Range("Q7").Value = 52.976913006825477
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = -0.0944
$ws.Range("T7").Value = -0.0257
$ws.Range("U7").Value = 6.62
$ws.Range("V7").Value = 1.88
$ws.Range("W7").Value = -2
